# python auto update tool
# Update the generated title-block layout IDs (sheet/layout tags) and
# restore the last-used selection on the Title Block sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "S1"
$ws.Range("A3").Value = "S2"
$ws.Range("A4").Value = "S3"

[void]$ws.Range("D20").Select()
